$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.16000000000049
$ws.Range("G2").Value = 0.01042232879617899
$ws.Range("H2").Value = 0.02365367902395644
$ws.Range("K2").Value = 4.946449009628239
$ws.Range("L2").Value = "[0.8051117765730673, 9.08778624268341]"
$ws.Range("M2").Value = 0.01944018520524526
$ws.Range("N2").Value = 0.01944018520524526
$ws.Range("O2").Value = -1.94973718240808
$ws.Range("P2").Value = "[-2.956053147521927, -0.9434212172942322]"
$ws.Range("Q2").Value = 0.0001725167627542046
$ws.Range("R2").Value = 0.0001725167627542046
$ws.Range("S2").Value = 15.29081846902388
$ws.Range("T2").Value = "[13.057168192999747, 17.52446874504801]"
$ws.Range("W2").Value = 7.80740740740756
$ws.Range("X2").Value = 3.777777777777852
$ws.Range("Y2").Value = 11.83703703703727

# Row 3 updates
$ws.Range("E3").Value = 25.03000000000047
$ws.Range("G3").Value = 0.0001199579243174043
$ws.Range("H3").Value = 0.001299863071920389
$ws.Range("K3").Value = 5.706677134143157
$ws.Range("L3").Value = "[2.1710696418498188, 9.242284626436495]"
$ws.Range("M3").Value = 0.001651620114569541
$ws.Range("N3").Value = 0.003303240229139082
$ws.Range("O3").Value = -2.025210879791619
$ws.Range("P3").Value = "[-2.6793162571156195, -1.3711055024676178]"
$ws.Range("Q3").Value = [double]"3.559691652554875e-09"
$ws.Range("R3").Value = [double]"7.11938330510975e-09"
$ws.Range("S3").Value = 13.19979006316723
$ws.Range("T3").Value = "[11.349279397261874, 15.050300729072589]"
$ws.Range("W3").Value = 8.067727727727881
$ws.Range("X3").Value = 5.462002002002107
$ws.Range("Y3").Value = 10.67345345345366
